$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$c = $ws.Cells.Item(200, 7)
$c.Value = "hello world"
Write-Host ("Value2=" + $c.Value2)
